$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lowRows  = @(3,6,9,12,15,18,21,24,27,30,33)
$highRows = @(4,7,10,13,16,19,22,25,28,31,34)

# ---------------------------------------------------------------------------
# 1. New small helper table to the right of the main table (I6:I7, J5:K6)
# ---------------------------------------------------------------------------
$ws.Range("I6").Value = "low"
$ws.Range("I7").Value = "high"

$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("J6").Value = 0.5
$ws.Range("K6").Value = 0

# ---------------------------------------------------------------------------
# 2. Update the row-label text in column B: "z = low"/"z = high" -> "z_t-1 = low"/"z_t-1 = high"
#    (style/formatting of these cells is unchanged)
# ---------------------------------------------------------------------------
foreach ($r in $lowRows) {
    $ws.Range("B$r").Value = "z_t-1 = low"
}
foreach ($r in $highRows) {
    $ws.Range("B$r").Value = "z_t-1 = high"
}

# ---------------------------------------------------------------------------
# 3. Add the new "z_t= low" / "z_t= high" labels in column E for every data
#    row, styled the same way as the column-B labels (bold, centered, General
#    number format).
# ---------------------------------------------------------------------------
foreach ($r in $lowRows) {
    $cell = $ws.Range("E$r")
    $cell.Value = "z_t= low"
    $cell.ClearFormats()
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}
foreach ($r in $highRows) {
    $cell = $ws.Range("E$r")
    $cell.Value = "z_t= high"
    $cell.ClearFormats()
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 4. Highlight cell C10 in bold (new emphasis style used by the author)
# ---------------------------------------------------------------------------
$ws.Range("C10").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Column sizing: column B now auto-sized to fit "z_t-1 = low/high", and
#    column E widened to comfortably fit the new "z_t= low/high" labels.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.083333333333334
$ws.Columns.Item(5).ColumnWidth = 12.25

# ---------------------------------------------------------------------------
# 6. Update the selected cell to reflect where the author ended up (M10)
# ---------------------------------------------------------------------------
$ws.Range("M10").Select()
